$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 and J1, matching style of existing header cells
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data for columns I (I0) and J (IF) for rows 2-25
$values = @{
    2  = @(8, 8)
    3  = @(8, 8)
    4  = @(7, 7)
    5  = @(7, 7)
    6  = @(5, 5)
    7  = @(6, 6)
    8  = @(7, 7)
    9  = @(6, 6)
    10 = @(8, 8)
    11 = @(8, 8)
    12 = @(6, 6)
    13 = @(3, 3)
    14 = @(6, 7)
    15 = @(8, 8)
    16 = @(8, 8)
    17 = @(7, 8)
    18 = @(6, 7)
    19 = @(7, 7)
    20 = @(9, 9)
    21 = @(8, 8)
    22 = @(8, 8)
    23 = @(6, 6)
    24 = @(6, 7)
    25 = @(7, 7)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
